$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (A: region, B: value, C: rank)
$ws.Range("A2").Value = "BR"
$ws.Range("B2").Value = 43.433333333333316

$ws.Range("A3").Value = "NE"
$ws.Range("B3").Value = 44.204444444444441

$ws.Range("A4").Value = "SE"
$ws.Range("B4").Value = 43.76
$ws.Range("C4").Value = 12

$ws.Range("A5").Value = "RJ"
$ws.Range("B5").Value = 47.47
$ws.Range("C5").Value = 6

$ws.Range("A6").Value = "PB"
$ws.Range("B6").Value = 47.88
$ws.Range("C6").Value = 5

$ws.Range("A7").Value = "AC"
$ws.Range("B7").Value = 48.41
$ws.Range("C7").Value = 4

$ws.Range("A8").Value = "RR"
$ws.Range("B8").Value = 49.28
$ws.Range("C8").Value = 3

$ws.Range("A9").Value = "MG"
$ws.Range("B9").Value = 51.37
$ws.Range("C9").Value = 2

$ws.Range("A10").Value = "RN"
$ws.Range("B10").Value = 56.94
$ws.Range("C10").Value = 1

# Update the selection to A2:C10 with active cell A2
$ws.Range("A2:C10").Select()
